$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column holds plain text values that look like ISO dates
# (e.g. "2026-02-19"). Force Text format on those two cells first so
# Excel stores the literal string instead of auto-converting it to a
# date serial number.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A9").NumberFormat = "@"

# Row 8: 2026-02-19
$ws.Range("A8").Value = "2026-02-19"
$ws.Range("B8").Value = "Booked"
$ws.Range("C8").Value = 12800
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "Auto-generated from bookings"

# Row 9: 2026-02-20
$ws.Range("A9").Value = "2026-02-20"
$ws.Range("B9").Value = "Booked"
$ws.Range("C9").Value = 12800
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "Auto-generated from bookings"
